$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first four data rows (rows 2-5), which corresponded to the
# evaluation years that should no longer be part of the series. This
# shifts all subsequent rows up by four and shrinks the used range from
# A1:B43 down to A1:B39.
$ws.Range("A2:B5").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
